# Refresh the cryptocurrency price/volume snapshot on sheet1 (data captured
# Thu Sep  7 16:29:02 UTC 2023). Every data cell on this sheet is stored as
# plain text (even the numeric-looking "Price" column), so for any new price
# string that Excel would otherwise auto-detect as a number we briefly apply a
# "@" (text) number format before writing the value, then clear that helper
# format again so the cell keeps the workbook default (no explicit) style -
# exactly like every other text cell already on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

Set-TextValue "D2" '25.754.25'
Set-TextValue "E2" '  +0.27%  '
Set-TextValue "D3" '1.628.82'
Set-TextValue "E3" '  -0.09%  '
Set-TextValue "D4" '0.998'
Set-TextValue "E4" '  -0.66%  '
Set-TextValue "D5" '214.44'
Set-TextValue "E5" '  +0.12%  '
Set-TextValue "E6" '  +0.04%  '
Set-TextValue "D7" '0.998'
Set-TextValue "E7" '  -0.66%  '
Set-TextValue "E8" '  -0.51%  '
Set-TextValue "D9" '0.0633'
Set-TextValue "E9" '  -0.46%  '
Set-TextValue "D10" '19.59'
Set-TextValue "E10" '  +0.80%  '
Set-TextValue "D11" '0.0794'
Set-TextValue "E11" '  +1.41%  '
Set-TextValue "D12" '4.26'
Set-TextValue "E12" '  +0.54%  '
Set-TextValue "D13" '1.851.94'
Set-TextValue "E13" '  -0.14%  '
Set-TextValue "D14" '1.604.62'
Set-TextValue "E14" '  -1.60%  '
Set-TextValue "D15" '0.554'
Set-TextValue "E15" '  +0.63%  '
Set-TextValue "D16" '0.0₃0762'
Set-TextValue "E16" '  -0.49%  '
Set-TextValue "D17" '62.81'
Set-TextValue "E17" '  -0.42%  '
Set-TextValue "D18" '25.646.27'
Set-TextValue "E18" '  -0.31%  '
Set-TextValue "D19" '0.998'
Set-TextValue "E19" '  -0.61%  '
Set-TextValue "D20" '4.45'
Set-TextValue "E20" '  +0.81%  '
Set-TextValue "D21" '191.15'
Set-TextValue "E21" '  -1.44%  '
Set-TextValue "D22" '9.92'
Set-TextValue "E22" '  +0.03%  '
Set-TextValue "D23" '6.27'
Set-TextValue "E23" '  +1.21%  '
Set-TextValue "D24" '0.998'
Set-TextValue "E24" '  -0.68%  '
Set-TextValue "E25" '  +1.57%  '
Set-TextValue "D26" '142.16'
Set-TextValue "E26" '  +1.54%  '
Set-TextValue "E27" '  +3.39%  '
Set-TextValue "D28" '6.84'
Set-TextValue "E28" '  +0.64%  '
Set-TextValue "D29" '15.49'
Set-TextValue "E29" '  +0.11%  '
Set-TextValue "E30" '  +0.19%  '
Set-TextValue "D31" '0.0494'
Set-TextValue "E31" '  +1.90%  '
Set-TextValue "D32" '3.33'
Set-TextValue "E32" '  -0.15%  '
Set-TextValue "D33" '3.23'
Set-TextValue "E33" '  -0.54%  '
Set-TextValue "D34" '1.59'
Set-TextValue "E34" '  +0.14%  '
Set-TextValue "E35" '  -0.14%  '
Set-TextValue "E36" '  +1.35%  '
Set-TextValue "D37" '1.139.24'
Set-TextValue "E37" '  +3.18%  '
Set-TextValue "D38" '2.50'
Set-TextValue "E38" '  -2.06%  '
Set-TextValue "E39" '  -0.28%  '
Set-TextValue "D40" '0.0156'
Set-TextValue "E40" '  +0.18%  '
Set-TextValue "E41" '  -0.66%  '
Set-TextValue "E42" '  -1.09%  '
Set-TextValue "D43" '5.58'
Set-TextValue "E43" '  +0.16%  '
Set-TextValue "D44" '100.70'
Set-TextValue "E44" '  +0.91%  '
Set-TextValue "D45" '0.803'
Set-TextValue "E45" '  +1.05%  '
Set-TextValue "D46" '1.762.08'
Set-TextValue "E46" '  -0.05%  '
Set-TextValue "D47" '55.23'
Set-TextValue "E47" '  +0.53%  '
Set-TextValue "E48" '  +1.69%  '
Set-TextValue "E49" '  +6.25%  '
Set-TextValue "D51" '2.32'
Set-TextValue "E51" '  -3.09%  '
